$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'Real Tomayapo'
$ws.Range("E2").Value = 'The Strongest'
$ws.Range("F2").Value = 1.14
$ws.Range("G2").Value = 690
$ws.Range("H2").Value = 110
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.01
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 1.31
$ws.Range("O2").Value = 1.55
$ws.Range("P2").Value = 1.01
$ws.Range("Q2").Value = 10.5
$ws.Range("R2").Value = 1.01
$ws.Range("S2").Value = 1.68
$ws.Range("T2").Value = 1.1
$ws.Range("U2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("AB2").Value = 1.5
$ws.Range("D3").Value = 'Bolivar'
$ws.Range("E3").Value = 'San Antonio Bulo Bulo'
$ws.Range("F3").Value = 1.01
$ws.Range("G3").Value = 1.01
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 1.66
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 500
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 2.12
$ws.Range("AO3").Value = 1000
$ws.Range("A4").Value = 'Honduras Liga Nacional'
$ws.Range("C4").Value = '22:00:00'
$ws.Range("D4").Value = 'Olancho'
$ws.Range("E4").Value = 'CD Marathon'
$ws.Range("F4").Value = 11.5
$ws.Range("G4").Value = 13
$ws.Range("H4").Value = 1.48
$ws.Range("I4").Value = 1.52
$ws.Range("J4").Value = 3.95
$ws.Range("K4").Value = 4.1
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 3.3
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 1.47
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 1.13
$ws.Range("S4").Value = 8.4
$ws.Range("T4").Value = 2.02
$ws.Range("U4").Value = 1.85
$ws.Range("V4").Value = 2.9
$ws.Range("W4").Value = 1.08
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 3.5
$ws.Range("Z4").Value = 6.4
$ws.Range("AA4").Value = 24
$ws.Range("AC4").Value = 5
$ws.Range("AD4").Value = 10.5
$ws.Range("AE4").Value = 55
$ws.Range("AG4").Value = 18.5
$ws.Range("AH4").Value = 40
$ws.Range("AI4").Value = 190
$ws.Range("AK4").Value = 130
$ws.Range("AL4").Value = 240
$ws.Range("AN4").Value = 620
$ws.Range("AO4").Value = 100
$ws.Range("F5").Value = 2.12
$ws.Range("G5").Value = 2.22
$ws.Range("H5").Value = 3.95
$ws.Range("I5").Value = 4.3
$ws.Range("J5").Value = 3.35
$ws.Range("K5").Value = 3.55
$ws.Range("L5").Value = 1.04
$ws.Range("N5").Value = 1.1
$ws.Range("O5").Value = 1.35
$ws.Range("P5").Value = 1.77
$ws.Range("Q5").Value = 1.3
$ws.Range("R5").Value = 1.36
$ws.Range("S5").Value = 3.6
$ws.Range("T5").Value = 1.05
$ws.Range("U5").Value = 1.04
$ws.Range("V5").Value = 1.01
$ws.Range("W5").Value = 1.02
$ws.Range("X5").Value = 18.5
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 990
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 16.5
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 970
$ws.Range("AI5").Value = 980
$ws.Range("AJ5").Value = 36
$ws.Range("AK5").Value = 980
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 990
